$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 and 40 swap (Algorand/Hedera) plus value updates
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06125"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2115"
$ws.Range("E40").Value = "  -3.07%  "

$ws.Range("D2").Value = "28.112.25"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.798.86"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "323.74"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -3.34%  "
$ws.Range("D8").Value = "0.3623"
$ws.Range("E8").Value = "  -3.19%  "
$ws.Range("D9").Value = "44.79"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "0.07537"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "1.121"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "21.63"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "6.194"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "7.363"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "1.818.15"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").Value = "92.75"
$ws.Range("E17").Value = "  +4.61%  "
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "0.06333"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "17.23"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "5.965"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "28.139.66"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "11.40"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").Value = "2.160"
$ws.Range("E25").Value = "  -6.85%  "
$ws.Range("D26").Value = "159.41"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").Value = "20.38"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "2.023.62"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("D29").Value = "2.228"
$ws.Range("E29").Value = "  -6.08%  "
$ws.Range("D30").Value = "127.96"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "1.173"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").Value = "5.894"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").Value = "0.09024"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "3.517"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("D35").Value = "12.82"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "0.02363"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").Value = "5.125"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Value = "0.6492"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D41").Value = "1.195"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "1.426"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "7.963"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "13.58"
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").Value = "0.6008"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "3.710"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "125.04"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "1.983"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "1.160"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").Value = "0.06972"
$ws.Range("E51").Value = "  +1.07%  "
